$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G ("MemoCol") - this shifts old G:J (FloatCol..CurrencyCol) to H:K
$ws.Columns("G").Insert()

# New column header
$ws.Range("G1").Value = "MemoCol"

# New memo data for rows 2 and 3 (string values)
$ws.Range("G2").Value = "Memo1"
$ws.Range("G3").Value = "Text 1"

# New "large int" numeric value for row 4, using the big-precision number format
$ws.Range("G4").NumberFormat = "0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000"
$ws.Range("G4").Value = 1.2

# Move the active selection to H6
$ws.Range("H6").Select()
